$wb = $excel.ActiveWorkbook

# Update the raw input value on Sheet1 (G3). This feeds the F3 formula
# (=F2+G3) which in turn feeds the chart's "Actual" series cache.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("G3").Value = 68.599999999999994

# Make Sheet1 the active/selected sheet (it becomes the visible tab).
$sheet1.Activate()
$sheet1.Range("G4").Select()
